# Apply the "added results of training demo for deeplab" edit.
#
# Summary of the change:
#  - On the "cloud" sheet, row 4 (the Deeplabv3 row) gets filled in with the
#    results of a training run: D4="10.2", E4="TODO", F4=14136, G4="TODO",
#    H4="TODO", I4="TODO", J4="TODO".
#  - The "cloud" sheet becomes the active/selected sheet (tabSelected moves
#    from "local" to "cloud"), with the cursor/selection left on G5.

$wb = $excel.ActiveWorkbook

$wsLocal = $wb.Worksheets.Item("local")
$wsCloud = $wb.Worksheets.Item("cloud")

# Fill in the new training-demo results on the "cloud" sheet, row 4.
# D4 holds a version-like label ("10.2") that must stay text, not be
# coerced into the number 10.2.
$wsCloud.Range("D4").NumberFormat = "@"
$wsCloud.Range("D4").Value = "10.2"
$wsCloud.Range("E4").Value = "TODO"
$wsCloud.Range("F4").Value = 14136
$wsCloud.Range("G4").Value = "TODO"
$wsCloud.Range("H4").Value = "TODO"
$wsCloud.Range("I4").Value = "TODO"
$wsCloud.Range("J4").Value = "TODO"

# Move the active sheet / selection to "cloud", with the cursor on G5.
$wsCloud.Activate()
$wsCloud.Range("G5").Select()
